$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 200
$ws.Range("I33").Value = 200
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 200
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 29
$ws.Range("N33").ClearContents()
$ws.Range("H58").Value = 2166.6667
$ws.Range("J58").Value = 2166.6667
$ws.Range("L58").Value = 6500.000100000001
$ws.Range("N58").Value = -6800.000100000001
$ws.Range("H76").Value = 3860.375
$ws.Range("I76").Value = 3697.7144
$ws.Range("J76").Value = 4999
$ws.Range("K76").Value = 3697.7144
$ws.Range("L76").Value = 4999
$ws.Range("M76").Value = -3382.7144
$ws.Range("N76").Value = -5629
$ws.Range("H79").Value = 3860.375
$ws.Range("I79").Value = 3697.7144
$ws.Range("J79").Value = 4999
$ws.Range("K79").Value = 3697.7144
$ws.Range("L79").Value = 4999
$ws.Range("M79").Value = -2605.7144
$ws.Range("N79").Value = -7183
$ws.Range("H94").Value = 4333.3335
$ws.Range("I94").Value = 4333.3335
$ws.Range("K94").Value = 4333.3335
$ws.Range("M94").Value = -3882.3335
$ws.Range("H113").Value = 3490.6365
$ws.Range("I113").Value = 3377.4443
$ws.Range("K113").Value = 3377.4443
$ws.Range("M113").Value = -123.4443000000001
$ws.Range("H116").Value = 13877.647
$ws.Range("I116").Value = 15953.625
$ws.Range("K116").Value = 15953.625
$ws.Range("M116").Value = -12511.625
$ws.Range("H132").Value = 41358.8
$ws.Range("I132").Value = 41358.8
$ws.Range("K132").Value = 124076.4
$ws.Range("M132").Value = -121546.4
$ws.Range("H136").Value = 91998.336
$ws.Range("J136").Value = 91998.336
$ws.Range("L136").Value = 91998.336
$ws.Range("N136").Value = -102198.336
$ws.Range("H137").Value = 17863104
$ws.Range("I137").Value = 22729474
$ws.Range("K137").Value = 68188422
$ws.Range("M137").Value = -68185872
$ws.Range("H138").Value = 5304.2173

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 1790
$ws.Range("I14").Value = 185
$ws.Range("K14").Value = 185
$ws.Range("M14").Value = -10
$ws.Range("H32").Value = 928811.5600000001
$ws.Range("I32").Value = 1011263.5
$ws.Range("K32").Value = 1011263.5
$ws.Range("M32").Value = -1010976.5
$ws.Range("H45").Value = 1456.1818
$ws.Range("I45").Value = 1045.5714
$ws.Range("J45").Value = 2174.75
$ws.Range("K45").Value = 1045.5714
$ws.Range("L45").Value = 2174.75
$ws.Range("M45").Value = -668.5714
$ws.Range("N45").Value = -2928.75
$ws.Range("H55").Value = 15362.5
$ws.Range("J55").Value = 29725
$ws.Range("L55").Value = 29725
$ws.Range("N55").Value = -30355
$ws.Range("H61").Value = 2859152
$ws.Range("I61").Value = 2068.1765
$ws.Range("K61").Value = 2068.1765
$ws.Range("M61").Value = -1856.1765
$ws.Range("H74").Value = 741743.5
$ws.Range("I74").Value = 808329.5
$ws.Range("J74").Value = 22614.8
$ws.Range("K74").Value = 808329.5
$ws.Range("L74").Value = 22614.8
$ws.Range("M74").Value = -807455.5
$ws.Range("N74").Value = -24362.8
$ws.Range("H77").Value = 741743.5
$ws.Range("I77").Value = 808329.5
$ws.Range("J77").Value = 22614.8
$ws.Range("K77").Value = 4041647.5
$ws.Range("L77").Value = 113074
$ws.Range("M77").Value = -4037279.5
$ws.Range("N77").Value = -121810
$ws.Range("H88").Value = 1946.0416
$ws.Range("I88").Value = 1060.6154
$ws.Range("J88").Value = 2992.4546
$ws.Range("K88").Value = 1060.6154
$ws.Range("L88").Value = 2992.4546
$ws.Range("M88").Value = -654.6153999999999
$ws.Range("N88").Value = -3804.4546
$ws.Range("H91").Value = 1946.0416
$ws.Range("I91").Value = 1060.6154
$ws.Range("J91").Value = 2992.4546
$ws.Range("K91").Value = 1060.6154
$ws.Range("L91").Value = 2992.4546
$ws.Range("M91").Value = 343.3846000000001
$ws.Range("N91").Value = -5800.4546
$ws.Range("H97").Value = 1356.8125
$ws.Range("I97").Value = 1356.8125
$ws.Range("K97").Value = 1356.8125
$ws.Range("M97").Value = -860.8125
$ws.Range("H110").Value = 2000
$ws.Range("I110").Value = 2000
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 2000
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 45
$ws.Range("N110").Value = -6090
$ws.Range("H132").Value = 4327.5264
$ws.Range("I132").Value = 2434.2727
$ws.Range("K132").Value = 7302.8181
$ws.Range("M132").Value = -4772.8181
$ws.Range("H136").Value = 2859152
$ws.Range("I136").Value = 2068.1765
$ws.Range("K136").Value = 6204.529500000001
$ws.Range("M136").Value = -3654.529500000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 23274.318
$ws.Range("I20").Value = 26359.705
$ws.Range("J20").Value = 12784
$ws.Range("K20").Value = 26359.705
$ws.Range("L20").Value = 12784
$ws.Range("M20").Value = -26112.705
$ws.Range("N20").Value = -13278
$ws.Range("H86").Value = 2959.1365
$ws.Range("I86").Value = 2242.4285
$ws.Range("K86").Value = 2242.4285
$ws.Range("M86").Value = -1119.4285
$ws.Range("H89").Value = 2959.1365
$ws.Range("I89").Value = 2242.4285
$ws.Range("K89").Value = 11212.1425
$ws.Range("M89").Value = -5596.1425
$ws.Range("H94").Value = 4662.2104
$ws.Range("I94").Value = 4034.2942
$ws.Range("J94").Value = 9999.5
$ws.Range("K94").Value = 4034.2942
$ws.Range("L94").Value = 9999.5
$ws.Range("M94").Value = -3583.2942
$ws.Range("N94").Value = -10901.5
$ws.Range("H99").Value = 17681.166
$ws.Range("I99").Value = 17681.166
$ws.Range("K99").Value = 17681.166
$ws.Range("M99").Value = -16183.166
$ws.Range("H105").Value = 6856.125
$ws.Range("I105").Value = 1890
$ws.Range("K105").Value = 1890
$ws.Range("M105").Value = -143
$ws.Range("H107").Value = 1514.7693
$ws.Range("I107").Value = 1210.5625
$ws.Range("J107").Value = 2001.5
$ws.Range("K107").Value = 1210.5625
$ws.Range("L107").Value = 2001.5
$ws.Range("M107").Value = 709.4375
$ws.Range("N107").Value = -5841.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 44194.39
$ws.Range("J86").Value = 10344.682
$ws.Range("L86").Value = 10344.682
$ws.Range("N86").Value = -12590.682
$ws.Range("H89").Value = 44194.39
$ws.Range("J89").Value = 10344.682
$ws.Range("L89").Value = 51723.41
$ws.Range("N89").Value = -62955.41
$ws.Range("H94").Value = 18816.334
$ws.Range("I94").Value = 27224.5
$ws.Range("K94").Value = 27224.5
$ws.Range("M94").Value = -26773.5
$ws.Range("H107").Value = 519.9643
$ws.Range("I107").Value = 411.16666
$ws.Range("J107").Value = 715.8
$ws.Range("K107").Value = 411.16666
$ws.Range("L107").Value = 715.8
$ws.Range("M107").Value = 1508.83334
$ws.Range("N107").Value = -4555.8
$ws.Range("H122").Value = 21761
$ws.Range("I122").Value = 2239.8572
$ws.Range("J122").Value = 44535.668
$ws.Range("K122").Value = 6719.571599999999
$ws.Range("L122").Value = 133607.004
$ws.Range("M122").Value = -4269.571599999999
$ws.Range("N122").Value = -138507.004

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 968.5
$ws.Range("J98").Value = 966.8
$ws.Range("L98").Value = 2900.4
$ws.Range("N98").Value = -5896.4
$ws.Range("H122").Value = 1010042.44
$ws.Range("I122").Value = 2305833
$ws.Range("K122").Value = 20752497
$ws.Range("M122").Value = -20750047
$ws.Range("H139").Value = 6989
$ws.Range("I139").Value = 3160.818
$ws.Range("J139").Value = 11200
$ws.Range("K139").Value = 9482.454000000002
$ws.Range("L139").Value = 33600
$ws.Range("M139").Value = -4342.454000000002
$ws.Range("N139").Value = -43880

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8377.388000000001
$ws.Range("I70").Value = 8487.450999999999
$ws.Range("J70").Value = 8187.8335
$ws.Range("K70").Value = 8487.450999999999
$ws.Range("L70").Value = 8187.8335
$ws.Range("M70").Value = -8217.450999999999
$ws.Range("N70").Value = -8727.833500000001
$ws.Range("H73").Value = 8377.388000000001
$ws.Range("I73").Value = 8487.450999999999
$ws.Range("J73").Value = 8187.8335
$ws.Range("K73").Value = 8487.450999999999
$ws.Range("L73").Value = 8187.8335
$ws.Range("M73").Value = -7551.450999999999
$ws.Range("N73").Value = -10059.8335
$ws.Range("H102").Value = 1411.625
$ws.Range("I102").Value = 1362
$ws.Range("J102").Value = 1626.6666
$ws.Range("K102").Value = 1362
$ws.Range("L102").Value = 1626.6666
$ws.Range("M102").Value = 260
$ws.Range("N102").Value = -4870.6666
$ws.Range("H107").Value = 4475
$ws.Range("J107").Value = 5000
$ws.Range("L107").Value = 5000
$ws.Range("N107").Value = -8840
$ws.Range("H122").Value = 2588.7932
$ws.Range("I122").Value = 2617.12
$ws.Range("J122").Value = 2411.75
$ws.Range("K122").Value = 7851.36
$ws.Range("L122").Value = 7235.25
$ws.Range("M122").Value = -5401.36
$ws.Range("N122").Value = -12135.25
$ws.Range("H132").Value = 13972.214
$ws.Range("I132").Value = 7592.625
$ws.Range("J132").Value = 52249.75
$ws.Range("K132").Value = 22777.875
$ws.Range("L132").Value = 156749.25
$ws.Range("M132").Value = -20247.875
$ws.Range("N132").Value = -161809.25
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1220
$ws.Range("I126").Value = 985.7143
$ws.Range("J126").Value = 1766.6666
$ws.Range("K126").Value = 2957.1429
$ws.Range("L126").Value = 5299.9998
$ws.Range("M126").Value = -487.1428999999998
$ws.Range("N126").Value = -10239.9998
$ws.Range("H136").Value = 5632008.5
$ws.Range("I136").Value = 2718948
$ws.Range("K136").Value = 8156844
$ws.Range("M136").Value = -8154294
